$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: clear phone_number (D2) entirely; mellicode (E2) becomes an empty
# text value (cell stays present but blank).
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = "'"

# Row 3: update phone_number (D3) and mellicode (E3) with new values.
# Leading zero is significant, so force text with a leading apostrophe,
# same as typing it into Excel.
$ws.Range("D3").Value = "'02938423984"
$ws.Range("E3").Value = "'0239482309"

# Remove row 4 entirely, shifting remaining rows up.
$ws.Rows.Item(4).Delete()
